$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Time Period labels in column A: shift quarter-end month from 09 to 12
$ws.Range("A2").Value = "2014-12"
$ws.Range("A3").Value = "2015-12"
$ws.Range("A4").Value = "2016-12"
$ws.Range("A5").Value = "2017-12"
$ws.Range("A6").Value = "2018-12"
$ws.Range("A7").Value = "2019-12"
$ws.Range("A8").Value = "2020-12"
$ws.Range("A9").Value = "2021-12"
$ws.Range("A10").Value = "2022-12"
$ws.Range("A11").Value = "2023-12"
$ws.Range("A12").Value = "2024-12"

$ws.Range("A13").Value = "2014-12"
$ws.Range("A14").Value = "2015-12"
$ws.Range("A15").Value = "2016-12"
$ws.Range("A16").Value = "2017-12"
$ws.Range("A17").Value = "2018-12"
$ws.Range("A18").Value = "2019-12"
$ws.Range("A19").Value = "2020-12"
$ws.Range("A20").Value = "2021-12"
$ws.Range("A21").Value = "2022-12"
$ws.Range("A22").Value = "2023-12"
$ws.Range("A23").Value = "2024-12"

# Update numeric data values in columns B:F for rows 2-25
$ws.Range("B2").Value = 15837972
$ws.Range("C2").Value = 227858
$ws.Range("D2").Value = 55198614
$ws.Range("E2").Value = 33993547
$ws.Range("F2").Value = 10686288
$ws.Range("B3").Value = 15997388
$ws.Range("C3").Value = 231370
$ws.Range("D3").Value = 55862465
$ws.Range("E3").Value = 34178087
$ws.Range("F3").Value = 10758850
$ws.Range("B4").Value = 16150496
$ws.Range("C4").Value = 234275
$ws.Range("D4").Value = 56532151
$ws.Range("E4").Value = 34370085
$ws.Range("F4").Value = 10829217
$ws.Range("B5").Value = 16295265
$ws.Range("C5").Value = 236703
$ws.Range("D5").Value = 57232615
$ws.Range("E5").Value = 34559337
$ws.Range("F5").Value = 10891399
$ws.Range("B6").Value = 16475801
$ws.Range("C6").Value = 239501
$ws.Range("D6").Value = 57977671
$ws.Range("E6").Value = 34766030
$ws.Range("F6").Value = 10971268
$ws.Range("B7").Value = 16658691
$ws.Range("C7").Value = 242324
$ws.Range("D7").Value = 58761104
$ws.Range("E7").Value = 35004427
$ws.Range("F7").Value = 11064101
$ws.Range("B8").Value = 16839617
$ws.Range("C8").Value = 245594
$ws.Range("D8").Value = 59499758
$ws.Range("E8").Value = 35205576
$ws.Range("F8").Value = 11142274
$ws.Range("B9").Value = 17042979
$ws.Range("C9").Value = 249416
$ws.Range("D9").Value = 60393769
$ws.Range("E9").Value = 35446022
$ws.Range("F9").Value = 11261191
$ws.Range("B10").Value = 17269164
$ws.Range("C10").Value = 254024
$ws.Range("D10").Value = 61238550
$ws.Range("E10").Value = 35754594
$ws.Range("F10").Value = 11400953
$ws.Range("B11").Value = 17461826
$ws.Range("C11").Value = 259603
$ws.Range("D11").Value = 62003759
$ws.Range("E11").Value = 35978294
$ws.Range("F11").Value = 11518590
$ws.Range("B12").Value = 17722036
$ws.Range("C12").Value = 264116
$ws.Range("D12").Value = 62933151
$ws.Range("E12").Value = 36360053
$ws.Range("F12").Value = 11717813
$ws.Range("B13").Value = 15837972
$ws.Range("C13").Value = 227858
$ws.Range("D13").Value = 55198614
$ws.Range("E13").Value = 33993547
$ws.Range("F13").Value = 10686288
$ws.Range("B14").Value = 15997388
$ws.Range("C14").Value = 231370
$ws.Range("D14").Value = 55862465
$ws.Range("E14").Value = 34178087
$ws.Range("F14").Value = 10758850
$ws.Range("B15").Value = 16150496
$ws.Range("C15").Value = 234275
$ws.Range("D15").Value = 56532151
$ws.Range("E15").Value = 34370085
$ws.Range("F15").Value = 10829217
$ws.Range("B16").Value = 16295265
$ws.Range("C16").Value = 236703
$ws.Range("D16").Value = 57232615
$ws.Range("E16").Value = 34559337
$ws.Range("F16").Value = 10891399
$ws.Range("B17").Value = 16475801
$ws.Range("C17").Value = 239501
$ws.Range("D17").Value = 57977671
$ws.Range("E17").Value = 34766030
$ws.Range("F17").Value = 10971268
$ws.Range("B18").Value = 16658691
$ws.Range("C18").Value = 242324
$ws.Range("D18").Value = 58761104
$ws.Range("E18").Value = 35004427
$ws.Range("F18").Value = 11064101
$ws.Range("B19").Value = 16839617
$ws.Range("C19").Value = 245594
$ws.Range("D19").Value = 59499758
$ws.Range("E19").Value = 35205576
$ws.Range("F19").Value = 11142274
$ws.Range("B20").Value = 17042979
$ws.Range("C20").Value = 249416
$ws.Range("D20").Value = 60393769
$ws.Range("E20").Value = 35446022
$ws.Range("F20").Value = 11261191
$ws.Range("B21").Value = 17269164
$ws.Range("C21").Value = 254024
$ws.Range("D21").Value = 61238550
$ws.Range("E21").Value = 35754594
$ws.Range("F21").Value = 11400953
$ws.Range("B22").Value = 17461826
$ws.Range("C22").Value = 259603
$ws.Range("D22").Value = 62003759
$ws.Range("E22").Value = 35978294
$ws.Range("F22").Value = 11518590
$ws.Range("B23").Value = 17722036
$ws.Range("C23").Value = 264116
$ws.Range("D23").Value = 62933151
$ws.Range("E23").Value = 36360053
$ws.Range("F23").Value = 11717813
$ws.Range("B24").Value = 1.01559122203156
$ws.Range("C24").Value = 1.23880742378265
$ws.Range("D24").Value = 1.25871901931235
$ws.Range("E24").Value = 0.587811009777879
$ws.Range("F24").Value = 0.697346427695706
$ws.Range("B25").Value = 1.24533347717766
$ws.Range("C25").Value = 1.7375697148786
$ws.Range("D25").Value = 1.38137252977639
$ws.Range("E25").Value = 0.762978330759675
$ws.Range("F25").Value = 1.15525912037069
